$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (Method/BaseURL/Endpoint/API_Key/Lat/Lon/City/HTTPStatus)
# is no longer needed -- delete it and let the data rows shift up.
$ws.Rows(1).Delete()

# The first data row (now row 1, previously row 2) expected a 200 but the
# fixture should assert 201 for this case.
$ws.Range("H1").Value = 201

# Update the remembered selection to match what was last selected.
$ws.Range("H2").Select()
